# Fix the "Jupyter Lab" product name to "JupyterLab" (no space) on the
# "Appendix B - Tools Used" slide.
#
# Replacing only the substring "Jupyter Lab" (leaving ", Version 4.1.2"
# untouched) causes PowerPoint to split the original single run into two
# runs at the edit boundary - "JupyterLab" and ", Version 4.1.2" - which
# matches how the author's edit landed in the saved OOXML.

$p = $ppt.ActivePresentation

$target = "Jupyter Lab"
$replacement = "JupyterLab"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)

        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $full = $tr.Text

            if ($full -ne $null -and $full.Contains($target)) {
                $charIndex = $full.IndexOf($target)

                while ($charIndex -ge 0) {
                    $start = $charIndex + 1
                    $len = $target.Length

                    $sub = $tr.Characters($start, $len)
                    $sub.Text = $replacement

                    $full = $tr.Text
                    $searchFrom = $charIndex + $replacement.Length
                    if ($searchFrom -lt $full.Length) {
                        $charIndex = $full.IndexOf($target, $searchFrom)
                    } else {
                        $charIndex = -1
                    }
                }
            }
        }
    }
}
